# The "Recorded By" column (G) lists the authors of each session record as a
# comma-separated string. A number of rows were exported with the literal
# "System" token first (e.g. "System, dnasr281@gmail.com"); this fixes the
# ordering so "System"/"system" is listed last instead, by reversing the
# comma-separated list - matching the convention used by the rest of the
# sheet (e.g. "dnasr281@gmail.com, System").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = $value -split ', '

    if ($parts.Count -gt 1 -and $parts[0].ToLower() -eq 'system') {
        $reversed = $parts[($parts.Count - 1)..0]
        $cell.Value = [string]::Join(', ', $reversed)
    }
}
